# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (fund holdings for the quarter), placed
#    right before the "总计" (grand total) sheet, using the same layout as
#    the existing quarterly sheets.
# 2. Add a corresponding row to the "总计" summary sheet and renumber the
#    existing index column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$beforeSheet = $wb.Worksheets.Item("总计")

# Clone the "2021-Q4" sheet (same column layout/styles we need) and drop it
# in right before "总计" - this keeps all of the sheetPr/outline/pageSetup
# and style formatting identical to its siblings.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($beforeSheet)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template had 3 data rows (rows 2-4); we only need 1, so drop the extras.
$newSheet.Range("A3:A4").EntireRow.Delete()

# Re-resolve a fresh, live reference to the "总计" sheet now that the new
# sheet has been inserted/renamed (a handle captured before the insert stays
# bound to its original position rather than following the renamed sheet).
$totalSheet = $wb.Worksheets.Item("总计")

# --- Fill in the 2022-Q1 fund-holding data row -----------------------------
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "398021"

$newSheet.Range("C2").Value = "中海能源策略混合"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "23.96"

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "90.88"

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "3.16"

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.7571"

$newSheet.Range("H2").Value = 7

# --- Update the "总计" (summary) sheet --------------------------------------
# Insert a new top data row for 2022-Q1 and renumber the existing index
# column (A) accordingly.
$totalSheet.Range("A2:D2").EntireRow.Insert()
$totalSheet.Range("A2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.76

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# Restore the originally active sheet/tab (inserting/copying a sheet makes it
# active by default, but the source workbook keeps "2021-Q2" as the selected
# tab).
$wb.Worksheets.Item("2021-Q2").Activate()
